# Update the "percent tonnage per depth range" table: rows 105-143 (A:J)
# are recalculated from a model run, shifting the depth-bucket split so the
# underground6/underground7 columns (I/J) pick up tonnage previously
# reported only in column I (and none in J).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newVals = @(
    0.08333333334,
    0.08333333332000002,
    0.08333333333999998,
    0.1071428572,
    0.107142857,
    0.1071428572000001,
    0.1071428571999999,
    0.1071428571999999,
    0.107142857,
    0.1071428572000001
)

for ($r = 105; $r -le 143; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $newVals[$c - 1]
    }
}
